# Apply the pagerank.docx edit:
#  - paragraph 2 (was empty) gets "This is an algorithm"
#  - paragraph 3 ("TFIDF") keeps its text but loses the _GoBack bookmark
#  - paragraph 4 (was empty) gets "And this is another algorithm"
#  - paragraph 5 ("DIRECT HIT") is unchanged
#  - a brand-new paragraph 6 is appended with
#    "And of course this is also algorithm, no surprise." and the
#    _GoBack bookmark is moved there, positioned right after the run
#    (a zero-length bookmark, not wrapping any text).

$d = $word.ActiveDocument

# --- Fill the two originally-empty paragraphs with their new text ---
$d.Paragraphs.Item(2).Range.Text = "This is an algorithm"
$d.Paragraphs.Item(4).Range.Text = "And this is another algorithm"

# --- Append a brand-new last paragraph after "DIRECT HIT" ---
$d.Paragraphs.Item(5).Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item(6)

# Write the target sentence plus one throw-away trailing character. The
# trailing character keeps the bookmark position we are about to create
# from landing exactly on the document's final offset, which this COM
# host mishandles for zero-length bookmarks (it resets to offset 0).
$finalText = "And of course this is also algorithm, no surprise."
$p6.Range.Text = $finalText + "X"

# --- Move the (single, Word-managed) _GoBack bookmark to sit right
#     after the new run, before the throw-away character ---
$p6 = $d.Paragraphs.Item(6)
$bmPos = $p6.Range.Start + $finalText.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Remove the throw-away trailing character, leaving the bookmark
#     collapsed immediately after the run's text ---
$p6 = $d.Paragraphs.Item(6)
$trailing = $d.Range($p6.Range.End - 2, $p6.Range.End - 1)
$trailing.Delete()
